$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows at 28:29 (existing rows 28-30 shift down to 30-32).
# Excel's native row-insert copies formatting from the row above, which already
# gives the new I28/I29 cells the s="7" style used throughout column I.
$ws.Rows("28:29").Insert()

# ---- Row 28: "sandflow" / 流沙 quest ----
$ws.Cells.Item(28,1).Value = 42010019          # A28 Id
$ws.Cells.Item(28,2).Value = "流沙"             # B28 Name
$ws.Cells.Item(28,3).Value = 1                 # C28 Type
$ws.Cells.Item(28,4).Value = 0                 # D28 Level
$ws.Cells.Item(28,5).Value = 2                 # E28 Danger
$ws.Cells.Item(28,6).Value = "sandflow"        # F28 Ename
$ws.Cells.Item(28,7).Value = "sandflow"        # G28 Figue
$ws.Cells.Item(28,8).Value = "sandflow"        # H28 Script
$ws.Cells.Item(28,23).Value = 100              # W28 PunishFood
$ws.Cells.Item(28,24).Value = 50               # X28 PunishHealth

# ---- Row 29: "swamp" / 沼泽 quest ----
$ws.Cells.Item(29,1).Value = 42010020          # A29 Id
$ws.Cells.Item(29,2).Value = "沼泽"             # B29 Name
$ws.Cells.Item(29,3).Value = 1                 # C29 Type
$ws.Cells.Item(29,4).Value = 0                 # D29 Level
$ws.Cells.Item(29,5).Value = 2                 # E29 Danger
$ws.Cells.Item(29,6).Value = "swamp"           # F29 Ename
$ws.Cells.Item(29,7).Value = "swamp"           # G29 Figue
$ws.Cells.Item(29,8).Value = "swamp"           # H29 Script
$ws.Cells.Item(29,23).Value = 50               # W29 PunishFood
$ws.Cells.Item(29,24).Value = 100              # X29 PunishHealth

# Column I ("TriggerMulti") holds the literal text "true" (not a boolean) in
# this sheet. Setting .Value to the string "true" would coerce to a native
# Excel boolean, so instead copy an existing text "true" cell's value through
# PasteSpecial(xlPasteValues), which preserves the shared-string text type
# and leaves the destination's existing style (s="7") untouched.
$ws.Cells.Item(27,9).Copy()
$ws.Cells.Item(28,9).PasteSpecial(-4163)
$ws.Cells.Item(29,9).PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Grow the worksheet table (ListObject) to cover the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:AC32"))

# Match the author's final selection.
$ws.Range("A29").Select()
